$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Training Dashboard")

# Use a scratch cell, far away from the used range, to hold the new
# "LAST UPDATE" date as plain text (not an Excel date serial). Writing the
# text directly into the target cells via Value2 causes Excel to
# auto-recognize the "dd-mmm-yyyy" pattern and convert it into a date
# serial number (changing both the stored value and the cell's number
# format/style). Instead we prepare the text once in a formatted-as-text
# scratch cell and then copy/paste-values it into each destination cell,
# which transfers only the text value and leaves the destination cell's
# existing style/number-format untouched.
$scratch = $ws.Cells.Item(1000, 1)
$scratch.NumberFormat = "@"
$scratch.Value2 = "04-Nov-2025"
$scratch.Copy()

# Rows 3 through 24: decrement "PERIOD TO EXPIRE" (column H) by 1 day,
# and set "LAST UPDATE" (column I) to 04-Nov-2025 (previously 03-Nov-2025).
for ($row = 3; $row -le 24; $row++) {
    $hCell = $ws.Cells.Item($row, 8)   # column H
    $hCell.Value2 = [double]$hCell.Value2 - 1

    $iCell = $ws.Cells.Item($row, 9)   # column I
    $iCell.PasteSpecial(-4163)          # xlPasteValues
}

# Clean up the scratch cell and clipboard marching ants / clipboard mode.
$scratch.Clear()
$excel.CutCopyMode = $false
